# "Fixed so many errors"
# - Append newly-coded business rows to Sheet1 (NAICS lookups, incl. "NOT FOUND")
# - Un-hide (veryHidden -> hidden) the helper sheet ESRI_MAPINFO_SHEET

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($rng, $val) {
    # Force text storage (t="s") even for digit-only strings like "722511"
    # by flipping to a text number format before the write, then clearing
    # the format back off so no stray style sticks to the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("A4").Value = "fakeco"
Set-TextValue $ws.Range("B4") "NOT FOUND"

$ws.Range("A5").Value = "green lantern lounge"
Set-TextValue $ws.Range("B5") "722511"

$ws.Range("A6").Value = "leo’s country oven"
Set-TextValue $ws.Range("B6") "NOT FOUND"

$ws.Range("A7").Value = "thang long"
Set-TextValue $ws.Range("B7") "722511"

$ws.Range("A8").Value = "galco"
Set-TextValue $ws.Range("B8") "443142"

# The ESRI_MAPINFO_SHEET helper tab was "very hidden" (not reachable from the
# Excel UI at all); relax it to a normal hidden sheet (unhide via right-click).
$ws2 = $wb.Worksheets.Item("ESRI_MAPINFO_SHEET")
$ws2.Visible = 0
